$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh: update D (Price) and E (Volume 1h) text cells.
# D-column values are stored as literal text (dotted thousands separators, fixed
# decimals). Excel's Range.Value setter auto-coerces numeric-looking strings to
# real numbers, which would both change the stored type and normalize formatting
# (e.g. "605.70" -> 605.7, "1.00" -> 1). To preserve the exact text, force the
# cell to Text format before assigning, then restore the default "Normal" style
# so no stray formatting is left behind (only the cell's content changes).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.033.53'
$ws.Range('D2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.519.24'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '605.70'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.34'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.518.80'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.87'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.98%  '
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000215'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.117.27'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '31.67'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.83%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.516.25'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.229.80'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.69'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +8.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.41'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.39'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '435.79'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.612'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '79.90'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.656.13'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -4.77%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.84'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.77%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.32'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -5.23%  '
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.60'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.18%  '
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.39'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.513.10'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.80'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.91%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.91'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.64%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.01'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0893'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '169.39'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.33%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.09'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -9.99%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.44'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.897'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '29.04'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.45%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '45.68'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.06%  '
$ws.Range('E48').Value = '  +1.12%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.47'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.33%  '
$ws.Range('E50').Value = '  -3.38%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.986'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.05%  '
